# Apply the Fri Aug 25 2023 cryptos list refresh (prices + 1h volume deltas)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value would otherwise be auto-parsed as a Number
# by Excel's type inference; force them back to Text so the stored value
# (and its exact displayed digits, e.g. trailing zeros) is preserved verbatim.
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D13',
    'D15',
    'D17',
    'D20',
    'D21',
    'D22',
    'D23',
    'D25',
    'D26',
    'D27',
    'D28',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D43',
    'D47',
    'D48',
    'D51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.107.21'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.659.89'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '217.48'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').Value = '0.5228'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '0.2633'
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').Value = '0.06276'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').Value = '20.59'
$ws.Range('E10').Value = '  -3.71%  '
$ws.Range('D11').Value = '0.07742'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.676.41'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.466'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '1.890.86'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '0.5443'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '0.0₅8092'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').Value = '64.80'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').Value = '26.138.50'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').Value = '4.579'
$ws.Range('E20').Value = '  -2.51%  '
$ws.Range('D21').Value = '191.27'
$ws.Range('D22').Value = '10.00'
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('D23').Value = '5.986'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').Value = '137.85'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').Value = '0.1238'
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('D27').Value = '7.243'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').Value = '16.16'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('D30').Value = '0.05922'
$ws.Range('E30').Value = '  -2.40%  '
$ws.Range('D31').Value = '1.277'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').Value = '3.524'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').Value = '3.252'
$ws.Range('E33').Value = '  -3.62%  '
$ws.Range('D34').Value = '1.563'
$ws.Range('E34').Value = '  -6.03%  '
$ws.Range('D35').Value = '0.9527'
$ws.Range('E35').Value = '  -3.79%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').Value = '2.769'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '0.5662'
$ws.Range('E38').Value = '  -5.31%  '
$ws.Range('D39').Value = '0.01593'
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('D40').Value = '5.896'
$ws.Range('E40').Value = '  -1.49%  '
$ws.Range('D41').Value = '0.8481'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').Value = '100.92'
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('D44').Value = '1.000.70'
$ws.Range('E44').Value = '  -7.16%  '
$ws.Range('D45').Value = '1.806.18'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '56.58'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').Value = '0.9999'
$ws.Range('E48').Value = '  -0.36%  '
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('E50').Value = '  +1.84%  '
$ws.Range('D51').Value = '0.05154'
$ws.Range('E51').Value = '  -0.80%  '
